$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2995
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 2995
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 2995
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(32, 14).Value = -3647

$ws.Cells.Item(55, 8).Value = 858.0625
$ws.Cells.Item(55, 9).Value = 756.3333
$ws.Cells.Item(55, 10).Value = 919.1
$ws.Cells.Item(55, 11).Value = 756.3333
$ws.Cells.Item(55, 12).Value = 919.1
$ws.Cells.Item(55, 13).Value = -542.3333
$ws.Cells.Item(55, 14).Value = -1347.1

$ws.Cells.Item(69, 8).Value = 14612.5
$ws.Cells.Item(69, 10).Value = 14612.5
$ws.Cells.Item(69, 12).Value = 43837.5
$ws.Cells.Item(69, 14).Value = -45585.5

$ws.Cells.Item(72, 8).Value = 14612.5
$ws.Cells.Item(72, 10).Value = 14612.5
$ws.Cells.Item(72, 12).Value = 131512.5
$ws.Cells.Item(72, 14).Value = -140248.5

$ws.Cells.Item(92, 8).Value = 57837.832
$ws.Cells.Item(92, 9).Value = 73634.42999999999
$ws.Cells.Item(92, 10).Value = 2549.75
$ws.Cells.Item(92, 11).Value = 73634.42999999999
$ws.Cells.Item(92, 12).Value = 2549.75
$ws.Cells.Item(92, 13).Value = -72386.42999999999
$ws.Cells.Item(92, 14).Value = -5045.75

$ws.Cells.Item(113, 8).Value = 14589.75
$ws.Cells.Item(113, 9).Value = 16411.5
$ws.Cells.Item(113, 11).Value = 16411.5
$ws.Cells.Item(113, 13).Value = -13157.5

$ws.Cells.Item(137, 8).Value = 1927.2273
$ws.Cells.Item(137, 9).Value = 1524.909
$ws.Cells.Item(137, 10).Value = 2329.5454
$ws.Cells.Item(137, 11).Value = 4574.727000000001
$ws.Cells.Item(137, 12).Value = 6988.6362
$ws.Cells.Item(137, 13).Value = -2024.727000000001
$ws.Cells.Item(137, 14).Value = -12088.6362

$ws.Cells.Item(138, 8).Value = 17745.785
$ws.Cells.Item(138, 9).Value = 23171.81
$ws.Cells.Item(138, 11).Value = 69515.43000000001
$ws.Cells.Item(138, 13).Value = -64375.43000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2785.5908
$ws.Cells.Item(74, 9).Value = 2564.7188
$ws.Cells.Item(74, 10).Value = 3374.5833
$ws.Cells.Item(74, 11).Value = 2564.7188
$ws.Cells.Item(74, 12).Value = 3374.5833
$ws.Cells.Item(74, 13).Value = -1690.7188
$ws.Cells.Item(74, 14).Value = -5122.5833

$ws.Cells.Item(77, 8).Value = 2785.5908
$ws.Cells.Item(77, 9).Value = 2564.7188
$ws.Cells.Item(77, 10).Value = 3374.5833
$ws.Cells.Item(77, 11).Value = 12823.594
$ws.Cells.Item(77, 12).Value = 16872.9165
$ws.Cells.Item(77, 13).Value = -8455.594000000001
$ws.Cells.Item(77, 14).Value = -25608.9165

$ws.Cells.Item(132, 8).Value = 24272.895
$ws.Cells.Item(132, 9).Value = 28534.691
$ws.Cells.Item(132, 10).Value = 3496.625
$ws.Cells.Item(132, 11).Value = 85604.073
$ws.Cells.Item(132, 12).Value = 10489.875
$ws.Cells.Item(132, 13).Value = -83074.073
$ws.Cells.Item(132, 14).Value = -15549.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 5852.1665
$ws.Cells.Item(99, 9).Value = 5175.6
$ws.Cells.Item(99, 11).Value = 5175.6
$ws.Cells.Item(99, 13).Value = -3677.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1052.5555
$ws.Cells.Item(22, 9).Value = 394.44446
$ws.Cells.Item(22, 11).Value = 394.44446
$ws.Cells.Item(22, 13).Value = -44.44445999999999

$ws.Cells.Item(99, 8).Value = 3698.875
$ws.Cells.Item(99, 9).Value = 3662.0908
$ws.Cells.Item(99, 10).Value = 3779.8
$ws.Cells.Item(99, 11).Value = 3662.0908
$ws.Cells.Item(99, 12).Value = 3779.8
$ws.Cells.Item(99, 13).Value = -2164.0908
$ws.Cells.Item(99, 14).Value = -6775.8

$ws.Cells.Item(107, 8).Value = 2091.2058
$ws.Cells.Item(107, 9).Value = 248.6
$ws.Cells.Item(107, 10).Value = 3545.8948
$ws.Cells.Item(107, 11).Value = 248.6
$ws.Cells.Item(107, 12).Value = 3545.8948
$ws.Cells.Item(107, 13).Value = 1671.4
$ws.Cells.Item(107, 14).Value = -7385.8948

$ws.Cells.Item(122, 8).Value = 2287.5
$ws.Cells.Item(122, 9).Value = 2005.5
$ws.Cells.Item(122, 11).Value = 6016.5
$ws.Cells.Item(122, 13).Value = -3566.5

$ws.Cells.Item(126, 8).Value = 3698.875
$ws.Cells.Item(126, 9).Value = 3662.0908
$ws.Cells.Item(126, 10).Value = 3779.8
$ws.Cells.Item(126, 11).Value = 10986.2724
$ws.Cells.Item(126, 12).Value = 11339.4
$ws.Cells.Item(126, 13).Value = -8516.2724
$ws.Cells.Item(126, 14).Value = -16279.4

$ws.Cells.Item(132, 8).Value = 2697.6
$ws.Cells.Item(132, 9).Value = 2697.6
$ws.Cells.Item(132, 11).Value = 8092.799999999999
$ws.Cells.Item(132, 13).Value = -5562.799999999999

$ws.Cells.Item(134, 8).Value = 77132.14
$ws.Cells.Item(134, 9).Value = 87444.664
$ws.Cells.Item(134, 11).Value = 262333.992
$ws.Cells.Item(134, 13).Value = -259798.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 3473.6924
$ws.Cells.Item(117, 10).Value = 3473.6924
$ws.Cells.Item(117, 12).Value = 10421.0772
$ws.Cells.Item(117, 14).Value = -17305.0772

$ws.Cells.Item(129, 8).Value = 4160.2
$ws.Cells.Item(129, 9).Value = 4530
$ws.Cells.Item(129, 10).Value = 4067.75
$ws.Cells.Item(129, 11).Value = 13590
$ws.Cells.Item(129, 12).Value = 12203.25
$ws.Cells.Item(129, 13).Value = -8590
$ws.Cells.Item(129, 14).Value = -22203.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4987.5557
$ws.Cells.Item(70, 9).Value = 4818.2
$ws.Cells.Item(70, 10).Value = 5199.25
$ws.Cells.Item(70, 11).Value = 4818.2
$ws.Cells.Item(70, 12).Value = 5199.25
$ws.Cells.Item(70, 13).Value = -4548.2
$ws.Cells.Item(70, 14).Value = -5739.25

$ws.Cells.Item(73, 8).Value = 4987.5557
$ws.Cells.Item(73, 9).Value = 4818.2
$ws.Cells.Item(73, 10).Value = 5199.25
$ws.Cells.Item(73, 11).Value = 4818.2
$ws.Cells.Item(73, 12).Value = 5199.25
$ws.Cells.Item(73, 13).Value = -3882.2
$ws.Cells.Item(73, 14).Value = -7071.25

$ws.Cells.Item(86, 8).Value = 59000
$ws.Cells.Item(86, 10).Value = 59000
$ws.Cells.Item(86, 12).Value = 59000
$ws.Cells.Item(86, 14).Value = -61372

$ws.Cells.Item(89, 8).Value = 59000
$ws.Cells.Item(89, 10).Value = 59000
$ws.Cells.Item(89, 12).Value = 177000
$ws.Cells.Item(89, 14).Value = -188856

$ws.Cells.Item(132, 8).Value = 37737.145
$ws.Cells.Item(132, 9).Value = 47363.273
$ws.Cells.Item(132, 10).Value = 2441.3333
$ws.Cells.Item(132, 11).Value = 142089.819
$ws.Cells.Item(132, 12).Value = 7323.999899999999
$ws.Cells.Item(132, 13).Value = -139559.819
$ws.Cells.Item(132, 14).Value = -12383.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2828.625
$ws.Cells.Item(7, 9).Value = 2173.7693
$ws.Cells.Item(7, 11).Value = 2173.7693
$ws.Cells.Item(7, 13).Value = -2061.7693

$ws.Cells.Item(22, 8).Value = 2095.3333
$ws.Cells.Item(22, 9).Value = 726
$ws.Cells.Item(22, 10).Value = 3036.75
$ws.Cells.Item(22, 11).Value = 726
$ws.Cells.Item(22, 12).Value = 3036.75
$ws.Cells.Item(22, 13).Value = -431
$ws.Cells.Item(22, 14).Value = -3626.75

$ws.Cells.Item(27, 8).Value = 2095.3333
$ws.Cells.Item(27, 9).Value = 726
$ws.Cells.Item(27, 10).Value = 3036.75
$ws.Cells.Item(27, 11).Value = 726
$ws.Cells.Item(27, 12).Value = 3036.75
$ws.Cells.Item(27, 13).Value = -619
$ws.Cells.Item(27, 14).Value = -3250.75

$ws.Cells.Item(46, 8).Value = 6432.6294
$ws.Cells.Item(46, 10).Value = 2549.2083
$ws.Cells.Item(46, 12).Value = 2549.2083
$ws.Cells.Item(46, 14).Value = -2925.2083

$ws.Cells.Item(61, 8).Value = 3780.348
$ws.Cells.Item(61, 9).Value = 3565.818
$ws.Cells.Item(61, 11).Value = 3565.818
$ws.Cells.Item(61, 13).Value = -3363.818

$ws.Cells.Item(113, 8).Value = 3780.348
$ws.Cells.Item(113, 9).Value = 3565.818
$ws.Cells.Item(113, 11).Value = 3565.818
$ws.Cells.Item(113, 13).Value = -1395.818

$ws.Cells.Item(126, 8).Value = 2828.625
$ws.Cells.Item(126, 9).Value = 2173.7693
$ws.Cells.Item(126, 11).Value = 6521.3079
$ws.Cells.Item(126, 13).Value = -4051.3079

$ws.Cells.Item(132, 8).Value = 37642.57
$ws.Cells.Item(132, 9).Value = 50232.36
$ws.Cells.Item(132, 11).Value = 150697.08
$ws.Cells.Item(132, 13).Value = -148167.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 2316.5
$ws.Cells.Item(113, 9).Value = 1899.75
$ws.Cells.Item(113, 10).Value = 3150
$ws.Cells.Item(113, 11).Value = 5699.25
$ws.Cells.Item(113, 12).Value = 9450
$ws.Cells.Item(113, 13).Value = -3529.25
$ws.Cells.Item(113, 14).Value = -13790

$ws.Cells.Item(122, 8).Value = 2247.5
$ws.Cells.Item(122, 9).Value = 2247.5
$ws.Cells.Item(122, 11).Value = 6742.5
$ws.Cells.Item(122, 13).Value = -4292.5
